$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $firstCol, $lastCol) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value() = $v2
        $cell2.Value() = $v1
    }
}

# Columns B..AC correspond to columns 2..29
Swap-Rows $ws 189 190 2 29
Swap-Rows $ws 192 193 2 29
Swap-Rows $ws 194 195 2 29
Swap-Rows $ws 226 227 2 29
